$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "November 10 - November 16; November 21 - November 24"
$ws.Range("E4").Value = "November 10 - November 16; November 21 - November 24"

$ws.Range("J4").Select()
